$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.3267545089723285
$ws.Range("C2").Value = -1.381391090384795
$ws.Range("D2").Value = 1.261504788818342
$ws.Range("E2").Value = 1.41980735745068
$ws.Range("F2").Value = 9.10735023142162
$ws.Range("B3").Value = 2.176221320692079
$ws.Range("C3").Value = 0.1327201135521534
$ws.Range("D3").Value = 9.137148433186667
$ws.Range("E3").Value = 11.12578741154633
$ws.Range("F3").Value = 2.054050900026041
$ws.Range("G3").Value = -2.583908390362894
$ws.Range("H3").Value = -1.69747595024657
$ws.Range("B4").Value = 5.072047478617549
$ws.Range("C4").Value = 7.687458084157399
$ws.Range("D4").Value = -0.3195551457329056
$ws.Range("E4").Value = -3.945495407503046
$ws.Range("F4").Value = -2.774906658561533
$ws.Range("B5").Value = -2.091464565649773
$ws.Range("C5").Value = -2.860361482790566
$ws.Range("D5").Value = -1.742507872682495
$ws.Range("E5").Value = -1.942182634720908
$ws.Range("F5").Value = -7.082866626106537
$ws.Range("G5").Value = -1.207897917465559
$ws.Range("H5").Value = -0.716788366044284
$ws.Range("B6").Value = 1.221728654477403
$ws.Range("C6").Value = -0.5859119406561888
$ws.Range("D6").Value = -6.412074126222691
$ws.Range("E6").Value = -1.327131657487215
$ws.Range("F6").Value = -1.040273951972324
$ws.Range("B7").Value = -3.911025385310296
$ws.Range("C7").Value = 0.1444249274195069
$ws.Range("D7").Value = -0.640273951972304
$ws.Range("E7").Value = -0.8163288124294192
$ws.Range("F7").Value = -1.055574095503602
$ws.Range("G7").Value = -0.7305900175239801
$ws.Range("H7").Value = 3.387154252625351
$ws.Range("B8").Value = -0.03107555164873999
$ws.Range("C8").Value = -1.322751909565398
$ws.Range("D8").Value = -1.555574095503573
$ws.Range("E8").Value = -1.130590017524014
$ws.Range("F8").Value = 3.087154252625439
$ws.Range("B9").Value = -0.457975156993101
$ws.Range("C9").Value = -0.8467238061668501
$ws.Range("D9").Value = 3.115058789684938
$ws.Range("E9").Value = -0.094982760252875
$ws.Range("F9").Value = -1.559991797019691
$ws.Range("G9").Value = 0.4693523401675459
$ws.Range("H9").Value = 4.335609787148272
$ws.Range("B10").Value = 1.287154252625385
$ws.Range("C10").Value = -1.726328558848905
$ws.Range("D10").Value = -2.575910395958857
$ws.Range("E10").Value = -0.0469927324174082
$ws.Range("F10").Value = 3.818840888438729
$ws.Range("B11").Value = -2.852896213663939
$ws.Range("C11").Value = -1.346992732417406
$ws.Range("D11").Value = 3.118840888438683
$ws.Range("E11").Value = 2.703774865812904
$ws.Range("F11").Value = -1.329268110997489
$ws.Range("G11").Value = 1.122467699931292
$ws.Range("H11").Value = -1.263707446020148
$ws.Range("B12").Value = 1.605332937026674
$ws.Range("C12").Value = 2.30459539371671
$ws.Range("D12").Value = -1.43016291878159
$ws.Range("E12").Value = 1.222790349166474
$ws.Range("F12").Value = -1.163823951986927
$ws.Range("B13").Value = -0.6367314425734349
$ws.Range("C13").Value = 0.9229439053107598
$ws.Range("D13").Value = -1.463935803707826
$ws.Range("E13").Value = 1.135800823836746
$ws.Range("F13").Value = 0.2460865825587319
$ws.Range("G13").Value = -0.4469970844674691
$ws.Range("H13").Value = 0.6021172494375269
$ws.Range("B14").Value = -0.763169241599825
$ws.Range("C14").Value = 1.301306938121286
$ws.Range("D14").Value = 0.332411999906014
$ws.Range("E14").Value = -0.363161970366406
$ws.Range("F14").Value = 0.6257449321420698
$ws.Range("B15").Value = -0.07192346926689197
$ws.Range("C15").Value = -0.593405114648007
$ws.Range("D15").Value = 0.3615079244737989
$ws.Range("E15").Value = 0.5676521307573097
$ws.Range("F15").Value = -0.4101892054037541
$ws.Range("G15").Value = 1.272006529183258
$ws.Range("H15").Value = 2.130826579980933
$ws.Range("B16").Value = -0.5943231312159014
$ws.Range("C16").Value = 0.1911117407590979
$ws.Range("D16").Value = -0.545820502981428
$ws.Range("E16").Value = 1.132498584751104
$ws.Range("F16").Value = 2.042107852101211
$ws.Range("B17").Value = -0.319302555566225
$ws.Range("C17").Value = 1.074786130482238
$ws.Range("D17").Value = 2.10356167974993
$ws.Range("E17").Value = 0.3923208492579415
$ws.Range("F17").Value = 0.03830692153388782
$ws.Range("G17").Value = 1.59724387475427
$ws.Range("H17").Value = -0.542996656587406
$ws.Range("B18").Value = 1.13405647582411
$ws.Range("C18").Value = 0.2906067856675916
$ws.Range("D18").Value = 0.1127521161105709
$ws.Range("E18").Value = 1.660991432238532
$ws.Range("F18").Value = -0.5156539280024219
$ws.Range("G18").Value = -0.1725649808533369
$ws.Range("H18").Value = 0.02392643965407593
$ws.Range("I18").Value = -0.6410204229136129
$ws.Range("J18").Value = -1.725647314577074
$ws.Range("B19").Value = -0.2949374036109162
$ws.Range("C19").Value = 1.585672713657559
$ws.Range("D19").Value = -0.499207999665203
$ws.Range("E19").Value = -0.18996464563375
$ws.Range("F19").Value = 0.05548488847716682
$ws.Range("G19").Value = -0.5600043896796199
$ws.Range("H19").Value = -1.636978867985007
$ws.Range("B20").Value = -0.6679969825334529
$ws.Range("C20").Value = -0.09868861252282102
$ws.Range("D20").Value = 0.03866155584581987
$ws.Range("E20").Value = -0.6015996449347509
$ws.Range("F20").Value = -1.6785706638791
$ws.Range("G20").Value = 2.099109144099572
$ws.Range("H20").Value = 0.3768551991671956
$ws.Range("I20").Value = 1.958424945524735
$ws.Range("J20").Value = 0.4058371456671066
$ws.Range("B21").Value = 0.7924721420310299
$ws.Range("C21").Value = -0.3251706130150041
$ws.Range("D21").Value = -1.653216592742552
$ws.Range("E21").Value = 2.108259859108429
$ws.Range("F21").Value = 0.3687611689862105
$ws.Range("G21").Value = 1.950332267839124
$ws.Range("H21").Value = 0.3977458174018127
$ws.Range("B22").Value = -1.114522562591048
$ws.Range("C22").Value = 2.578076292393661
$ws.Range("D22").Value = 0.6085133680368104
$ws.Range("E22").Value = 2.007807608059701
$ws.Range("F22").Value = 0.4552400051582457
$ws.Range("G22").Value = 0.1075696162591699
$ws.Range("H22").Value = 2.454548595455236
$ws.Range("I22").Value = 0.06585219932875908
$ws.Range("J22").Value = 1.341363851463769
$ws.Range("B23").Value = 2.343805781768054
$ws.Range("C23").Value = 0.5358128347403555
$ws.Range("D23").Value = 1.952179906974288
$ws.Range("E23").Value = 0.3997166069887186
$ws.Range("F23").Value = 0.05295683044919008
$ws.Range("G23").Value = 2.399557791075864
$ws.Range("H23").Value = 0.01096685068290504
$ws.Range("I23").Value = 1.286584226738754
$ws.Range("B24").Value = -0.2882865247592008
$ws.Range("C24").Value = 1.95265349408306
$ws.Range("D24").Value = 0.4001911081535217
$ws.Range("E24").Value = 0.053958022313495
$ws.Range("F24").Value = 2.400142854569665
$ws.Range("G24").Value = 0.01155303567007593
$ws.Range("H24").Value = 1.287171533861382
$ws.Range("B25").Value = 2.281065482681691
$ws.Range("C25").Value = 0.7301802256994567
$ws.Range("D25").Value = 0.05170297094063891
$ws.Range("E25").Value = 2.398961748446055
$ws.Range("F25").Value = -0.08963034146584925
$ws.Range("G25").Value = 1.185986011406541
$ws.Range("H25").Value = 3.852780295134727
$ws.Range("I25").Value = 21.08909300779241
$ws.Range("B26").Value = 0.04740873415159963
$ws.Range("C26").Value = 0.2
$ws.Range("D26").Value = 2.446768218203209
$ws.Range("E26").Value = -0.04195186392760797
$ws.Range("F26").Value = 1.233536122258116
$ws.Range("G26").Value = 3.9
$ws.Range("H26").Value = 21.13630385521235
$ws.Range("B27").Value = 0.712723949918967
$ws.Range("C27").Value = 2.564738550277127
$ws.Range("D27").Value = -0.08799717683646713
$ws.Range("E27").Value = 1.187622307891402
$ws.Range("F27").Value = 3.854858078752349
$ws.Range("G27").Value = 21.09093088973817
$ws.Range("B28").Value = 0.446768218203209
$ws.Range("C28").Value = 0.05804813607239212
$ws.Range("D28").Value = 1.233536122258116
$ws.Range("E28").Value = 3.9
$ws.Range("F28").Value = 21.13630385521235
$ws.Range("B29").Value = -0.12716270916539
$ws.Range("C29").Value = 1.271189067131602
$ws.Range("D29").Value = 3.894528479370229
$ws.Range("E29").Value = 21.13104881080667
$ws.Range("B30").Value = 1.385915977688484
$ws.Range("C30").Value = 3.892970719425895
$ws.Range("D30").Value = 21.13025136750149
$ws.Range("B31").Value = 1.541804305759959
$ws.Range("C31").Value = 17.97445767672242
$ws.Range("B32").Value = -4.063696144787649
